# Rename the existing sheet and add the new "pm10_limits" sheet right after it
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "arima_graph"

$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "pm10_limits"

# Fill in the PM10 limits comparison table (EU vs WHO/OMS)
$ws2.Range("F4").Value = "ue"
$ws2.Range("G4").Value = "oms"
$ws2.Range("E5").Value = "Promedio 24 horas"
$ws2.Range("F5").Value = "<= 50 ug/m3; <= 35 días"
$ws2.Range("G5").Value = "<= 50 ug/m3; <= 3 días"
$ws2.Range("E6").Value = "Promedio anual"
$ws2.Range("F6").Value = "<= 40 ug/m3"
$ws2.Range("G6").Value = "<= 20 ug/m3"

# Center-align the limit values
$ws2.Range("F4:G6").HorizontalAlignment = -4108

# Update the chart's series formulas to point at the renamed sheet
$cs = $ws1.ChartObjects().Item(1)
$chart = $cs.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(arima_graph!`$B`$1,,arima_graph!`$B`$2:`$B`$13,1)"
$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(arima_graph!`$C`$1,,arima_graph!`$C`$2:`$C`$13,2)"

# Make pm10_limits the active/selected sheet, with G9 as the last selection
$ws2.Activate()
$ws2.Range("G9").Select()
